$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2-11) from 45175 to 45183
$ws.Range("C2:C11").Value = 45183
